$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 11.347161
$ws.Range("H2").Value = 34.041483
$ws.Range("I2").Value = 0.1617975773769501
$ws.Range("J2").Value = 0.1617975773769501
$ws.Range("M2").Value = 0.5273163333333333
$ws.Range("N2").Value = 1.581949
$ws.Range("Q2").Value = 5.983543332262999
$ws.Range("R2").Value = 53.85188999036699
$ws.Range("S2").Value = 0.1617975773769501
$ws.Range("T2").Value = 0.1617975773769501

# Row 3
$ws.Range("I3").Value = 0.5781294662164954
$ws.Range("J3").Value = 0.5781294662164954
$ws.Range("M3").Value = 0.5273163333333333
$ws.Range("N3").Value = 1.581949
$ws.Range("Q3").Value = 21.38018855934544
$ws.Range("R3").Value = 192.421697034109
$ws.Range("S3").Value = 0.5781294662164954
$ws.Range("T3").Value = 0.5781294662164954

# Row 4
$ws.Range("G4").Value = 18.23939366666667
$ws.Range("H4").Value = 54.718181
$ws.Range("I4").Value = 0.2600729564065544
$ws.Range("J4").Value = 0.2600729564065544
$ws.Range("M4").Value = 0.5273163333333333
$ws.Range("N4").Value = 1.581949
$ws.Range("Q4").Value = 9.617930190529888
$ws.Range("R4").Value = 86.561371714769
$ws.Range("S4").Value = 0.2600729564065544
$ws.Range("T4").Value = 0.2600729564065544
